{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// ---------------------------------------------------------------------\n// 1. \"* as Export (Namespace Import)\" heading: merge the \"as\" run and the\n//    \" Export (Namespace Import)\" run (previously split apart by a\n//    grammar-check proofing mark) into a single \"as Export (Namespace\n//    Import)\" run.\n// ---------------------------------------------------------------------\nconst asRanges = context.document.body.search(\"as Export (Namespace Import)\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nasRanges.load(\"text\");\nawait context.sync();\n\nif (asRanges.items.length > 0) {\n  asRanges.items[0].insertText(\"as Export (Namespace Import)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2. Insert a new sentence at the very beginning of the paragraph that\n//    starts with \"The useState hook in React is essential...\".\n// ---------------------------------------------------------------------\nconst targetRanges = context.document.body.search(\"The useState hook in React is essential\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ntargetRanges.load(\"text\");\nawait context.sync();\n\nif (targetRanges.items.length > 0) {\n  const para = targetRanges.items[0].paragraphs.getFirst();\n  para.insertText(\n    \"Whenever a state variable changes, react re-renders the component. \",\n    Word.InsertLocation.start\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d ($word.ActiveDocument) are pre-seeded.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1. \"* as Export (Namespace Import)\" heading: merge the \"as\" run and the\n#    \" Export (Namespace Import)\" run into a single \"as Export (Namespace\n#    Import)\" run (dropping the grammar-check proofing marks around \"as\").\n# ---------------------------------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.Execute(\"as Export (Namespace Import)\", $false, $false, $false, $false, $false, $true, 1, $false, \"as Export (Namespace Import)\", 2)\n\n# ---------------------------------------------------------------------\n# 2. Insert a new sentence at the very beginning of the paragraph that\n#    starts with \"The useState hook in React is essential...\".\n# ---------------------------------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.Execute(\"The useState hook in React is essential\")\nif ($rng2.Find.Found) {\n    $para = $rng2.Paragraphs(1)\n    $insertPoint = $para.Range\n    $insertPoint.Collapse(1)\n    $insertPoint.InsertBefore(\"Whenever a state variable changes, react re-renders the component. \")\n}\n"}
